$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the AutoFilter criteria: un-hides all filtered-out rows and drops the
# stored filter value, while keeping the AutoFilter dropdown arrows/range.
$ws.ShowAllData()

# Widen column C to its new manually-set width (was a "best fit" width before).
$ws.Columns.Item(3).ColumnWidth = 22.28515625
